# Update the standard-error cells in table5_panel2 (bac_test_primary / multiple_imputation
# columns) to reflect the new pickled replication results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "(0.01)"
$ws.Range("C4").Value = "(0.12)"
$ws.Range("B6").Value = "(0.07)"
$ws.Range("C6").Value = "(0.12)"
